$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows that fall outside the new 4-row matrix (rows 5-13)
$ws.Range("A5:E13").Clear()

# Give A1 the same header formatting already used by B1:E1 (bold, bordered, centered)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row: Year + one column per region
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Eastern"
$ws.Range("C1").Value = "Northern"
$ws.Range("D1").Value = "Southern"
$ws.Range("E1").Value = "Westerns"

# Row 2 - 1996
$ws.Range("A2").Value = 1996
$ws.Range("B2").Value = "Ernst Handel"
$ws.Range("C2").Value = "Save-a-lot Markets"
$ws.Range("D2").Value = "LILA-Supermercado"
$ws.Range("E2").Value = "Piccolo und mehr"

# Row 3 - 1997
$ws.Range("A3").Value = 1997
$ws.Range("B3").Value = "QUICK-Stop"
$ws.Range("C3").Value = "QUICK-Stop"
$ws.Range("D3").Value = "QUICK-Stop"
$ws.Range("E3").Value = "Save-a-lot Markets"

# Row 4 - 1998
$ws.Range("A4").Value = 1998
$ws.Range("B4").Value = "QUICK-Stop"
$ws.Range("C4").Value = "Rattlesnake Canyon Grocery"
$ws.Range("D4").Value = "Hungry Owl All-Night Grocers"
$ws.Range("E4").Value = "Save-a-lot Markets"
